$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 3226
$ws.Range("K3").Value = 3192
$ws.Range("C4").Value = 1849
$ws.Range("K4").Value = 658
$ws.Range("K5").Value = 211
$ws.Range("K6").Value = 3786
$ws.Range("C7").Value = 28393
$ws.Range("K7").Value = 11073

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 214
$ws.Range("K4").Value = 41
$ws.Range("K6").Value = 239
$ws.Range("K7").Value = 737

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 74
$ws.Range("K6").Value = 58
$ws.Range("K7").Value = 239

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 122
$ws.Range("K3").Value = 158
$ws.Range("K6").Value = 124
$ws.Range("K7").Value = 435

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K2").Value = 59
$ws.Range("K3").Value = 65
$ws.Range("K7").Value = 182

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K3").Value = 132
$ws.Range("K7").Value = 383

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 82
$ws.Range("K3").Value = 67
$ws.Range("K4").Value = 8
$ws.Range("K6").Value = 102
$ws.Range("K7").Value = 265

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Range("K3").Value = 72
$ws.Range("K7").Value = 194

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 87
$ws.Range("K6").Value = 86
$ws.Range("K7").Value = 313
$ws.Range("K8").Value = 737
$ws.Range("K9").Value = 43
$ws.Range("K11").Value = 237
$ws.Range("K15").Value = 114
$ws.Range("K16").Value = 36
$ws.Range("K19").Value = 337
$ws.Range("K20").Value = 254
$ws.Range("K27").Value = 113
$ws.Range("K29").Value = 581
$ws.Range("K31").Value = 121
$ws.Range("K33").Value = 435
$ws.Range("K34").Value = 53
$ws.Range("K36").Value = 133
$ws.Range("K37").Value = 383
$ws.Range("K40").Value = 27
$ws.Range("K42").Value = 385
$ws.Range("K43").Value = 98
$ws.Range("K45").Value = 11
$ws.Range("K47").Value = 58
$ws.Range("K48").Value = 135
$ws.Range("K50").Value = 65
$ws.Range("K51").Value = 128
$ws.Range("K52").Value = 301
$ws.Range("K54").Value = 213
$ws.Range("K55").Value = 119
$ws.Range("K59").Value = 22
$ws.Range("K60").Value = 67
$ws.Range("C63").Value = 277
$ws.Range("K63").Value = 39
$ws.Range("K65").Value = 265
$ws.Range("K67").Value = 435
$ws.Range("K76").Value = 168
$ws.Range("K77").Value = 79
$ws.Range("K78").Value = 141
$ws.Range("K79").Value = 288
$ws.Range("K80").Value = 38
$ws.Range("K83").Value = 239
$ws.Range("K84").Value = 76
$ws.Range("K85").Value = 523
$ws.Range("K88").Value = 132
$ws.Range("K89").Value = 148
$ws.Range("K90").Value = 98
$ws.Range("K95").Value = 182
$ws.Range("K98").Value = 61
$ws.Range("K99").Value = 194
$ws.Range("K100").Value = 17
$ws.Range("C101").Value = 28393
$ws.Range("K101").Value = 11073

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 121

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 136
$ws.Range("K6").Value = 125
$ws.Range("K7").Value = 435

$ws = $wb.Worksheets.Item("South Deering")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 76

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K6").Value = 98
$ws.Range("K7").Value = 213

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K2").Value = 160
$ws.Range("K3").Value = 201
$ws.Range("K6").Value = 176
$ws.Range("K7").Value = 581

$ws = $wb.Worksheets.Item("Lake View")
$ws.Range("K6").Value = 71
$ws.Range("K7").Value = 135

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K6").Value = 110
$ws.Range("K7").Value = 337

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 31
$ws.Range("K6").Value = 99
$ws.Range("K7").Value = 168

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K3").Value = 120
$ws.Range("K6").Value = 150
$ws.Range("K7").Value = 385

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Range("K3").Value = 31
$ws.Range("K7").Value = 141

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K2").Value = 38
$ws.Range("K7").Value = 119

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K2").Value = 100
$ws.Range("K4").Value = 15
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 288

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Range("K2").Value = 87
$ws.Range("K4").Value = 9
$ws.Range("K7").Value = 254

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 53
$ws.Range("K4").Value = 14
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K6").Value = 9
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K2").Value = 107
$ws.Range("K7").Value = 313

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K3").Value = 17
$ws.Range("K7").Value = 53

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Range("K3").Value = 21
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K2").Value = 41
$ws.Range("K7").Value = 114

$ws = $wb.Worksheets.Item("Wicker Park")
$ws.Range("K6").Value = 45
$ws.Range("K7").Value = 61

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K3").Value = 8
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Range("K2").Value = 70
$ws.Range("K6").Value = 91
$ws.Range("K7").Value = 237

$ws = $wb.Worksheets.Item("Avalon Park")
$ws.Range("K3").Value = 16
$ws.Range("K7").Value = 43

$ws = $wb.Worksheets.Item("Montclare")
$ws.Range("K2").Value = 6
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 24
$ws.Range("K7").Value = 87

$ws = $wb.Worksheets.Item("United Center")
$ws.Range("K3").Value = 34
$ws.Range("K7").Value = 132

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K2").Value = 35
$ws.Range("K4").Value = 22
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 148

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K6").Value = 44
$ws.Range("K7").Value = 113

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Range("K4").Value = 8
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K3").Value = 32
$ws.Range("K7").Value = 128

$ws = $wb.Worksheets.Item("Morgan Park")
$ws.Range("K4").Value = 7
$ws.Range("K7").Value = 67

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 20
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 98

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 180
$ws.Range("K6").Value = 119
$ws.Range("K7").Value = 523

$ws = $wb.Worksheets.Item("Riverdale")
$ws.Range("K3").Value = 28
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Range("K5").Value = 1
$ws.Range("K7").Value = 11

$ws = $wb.Worksheets.Item("Rush & Division")
$ws.Range("K3").Value = 6
$ws.Range("K6").Value = 21
$ws.Range("K7").Value = 38

$ws = $wb.Worksheets.Item("Hegewisch")
$ws.Range("K6").Value = 5
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K2").Value = 79
$ws.Range("K3").Value = 79
$ws.Range("K7").Value = 301

$ws = $wb.Worksheets.Item("Bucktown")
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 36
